$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column I - match the formatting of the other header cells (H1)
$ws.Range("I1").Value = "eta²"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats

# New effect-size values for I2:I9
$values = @(0.04, 0.18, 0.08, 0.16, 0.19, 0.15, 0.27, 0.06)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}
